# PPT & PDF Update
# Applies the FSharp.pptx diff:
#  1. Slide 2 ("About F#"): merge the ". " + "JavaScript" runs into ". JavaScript"
#  2. Slide 4 ("History"):
#       - "From engineering side" -> "From engineering side (Turing-machine)"
#       - "From scientific side"  -> "From scientific side (\u03bb-calculus)"
#       - remove the leading "\u03bb-calculus " that used to prefix the
#         Wingdings arrow-chain paragraph, and split the merged
#         "IPL / LISP" Wingdings run into two separate runs
#  3. Slide 6 ("Don Syme's recommendation"): merge the " " + "(" runs into " ("

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 2 - "Can be comiled to i.e. JavaScript"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$para = $tr2.Paragraphs(5, 1)
$full = $para.Text
$idx = $full.IndexOf(". JavaScript")
$merge = $tr2.Characters($para.Start + $idx, 12)
$merge.Text = ". JavaScript"

# ---------------------------------------------------------------------------
# 2. Slide 4 - "History"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange

# 2a. "From engineering side" -> add " (Turing-machine)"
$paraEng = $tr4.Paragraphs(1, 1)
$afterEng = $paraEng.InsertAfter(" (Turing-machine)")
$fullEng = $paraEng.Text
$idxOpen = $fullEng.IndexOf(" (Turing-machine)")
$engBase = $paraEng.Start + $idxOpen
# split " (Turing-machine)" into " (" / "Turing-machine" / ")"
$r1 = $tr4.Characters($engBase, 2)
$r1.Text = " ("
$r2 = $tr4.Characters($engBase + 2, 14)
$r2.Text = "Turing-machine"
$r3 = $tr4.Characters($engBase + 16, 1)
$r3.Text = ")"

# 2b. "From scientific side" -> add " (\u03bb-calculus)"
$lambda = [string][char]0x03BB
$paraSci = $tr4.Paragraphs(5, 1)
$afterSci = $paraSci.InsertAfter(" (" + $lambda + "-calculus)")
$fullSci = $paraSci.Text
$idxOpen2 = $fullSci.IndexOf(" (" + $lambda + "-calculus)")
$sciBase = $paraSci.Start + $idxOpen2
# split into " (" / "\u03bb" / "-calculus" / ")"
$s1 = $tr4.Characters($sciBase, 2)
$s1.Text = " ("
$s2b = $tr4.Characters($sciBase + 2, 1)
$s2b.Text = $lambda
$s3 = $tr4.Characters($sciBase + 3, 9)
$s3.Text = "-calculus"
$s4c = $tr4.Characters($sciBase + 12, 1)
$s4c.Text = ")"

# 2c. Remove the leading "\u03bb-calculus " text that used to start the
#     "IPL -> LISP -> ..." bullet line (now redundant, see 2b above), and
#     split the merged Wingdings "IPL / LISP" run into two runs.
$paraList = $tr4.Paragraphs(7, 1)
$fullList = $paraList.Text
$idxLambda = $fullList.IndexOf($lambda)
$prefixRange = $tr4.Characters($paraList.Start + $idxLambda, 11)
$prefixRange.Text = ""

# The Wingdings run "<arrow> IPL <arrow> LISP <arrow>" now sits right after
# the leading tab character; replace/split it into "IPL " + " LISP ".
$wingStart = $paraList.Start + 1
$wingFull = $tr4.Characters($wingStart, 14)
$wingFull.Text = "IPL  LISP "
$wing1 = $tr4.Characters($wingStart, 4)
$wing1.Text = "IPL "
$wing2 = $tr4.Characters($wingStart + 4, 6)
$wing2.Text = " LISP "

# ---------------------------------------------------------------------------
# 3. Slide 6 - "If you don't need mutable values (state), use FP"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$para6 = $tr6.Paragraphs(2, 1)
$full6 = $para6.Text
$idx6 = $full6.IndexOf(" (state")
$merge6 = $tr6.Characters($para6.Start + $idx6, 2)
$merge6.Text = " ("
